$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.978.76'
$ws.Range("E2").Value = '  +1.51%  '
$ws.Range("D3").Value = '3.314.31'
$ws.Range("E3").Value = '  +6.08%  '
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '599.07'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.12%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '143.50'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +5.29%  '
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("D8").Value = '3.311.10'
$ws.Range("E8").Value = '  +6.23%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.524'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.48%  '
$ws.Range("E10").Value = '  +3.08%  '
$ws.Range("E11").Value = '  +5.80%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.475'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +4.13%  '
$ws.Range("E13").Value = '  +1.52%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '34.83'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.00%  '
$ws.Range("D15").Value = '3.859.49'
$ws.Range("E15").Value = '  +6.13%  '
$ws.Range("E16").Value = '  +1.28%  '
$ws.Range("D17").Value = '3.310.32'
$ws.Range("E17").Value = '  +6.17%  '
$ws.Range("D18").Value = '64.035.78'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.92'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +3.48%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '483.61'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.99%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.29'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.51%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.742'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +6.19%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.03'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +3.99%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '13.59'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +4.41%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '84.79'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -2.41%  '
$ws.Range("E26").Value = '  +0.25%  '
$ws.Range("E27").Value = '  +2.47%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.29'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.93%  '
$ws.Range("E29").Value = '  -0.09%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '8.23'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +3.93%  '
$ws.Range("B31").Value = 'EthereumClassic'
$ws.Range("C31").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '29.47'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +10.56%  '
$ws.Range("B32").Value = 'ImmutableX'
$ws.Range("C32").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.16'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +5.56%  '
$ws.Range("E33").Value = '  -0.34%  '
$ws.Range("E34").Value = '  +1.43%  '
$ws.Range("E35").Value = '  +2.84%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.01'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +3.40%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '53.32'
$ws.Range("D37").Style = "Normal"
$ws.Range("E38").Value = '  +8.23%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0401'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +3.46%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '432.59'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +3.03%  '
$ws.Range("E41").Value = '  +5.67%  '
$ws.Range("E42").Value = '  +2.60%  '
$ws.Range("E43").Value = '  +3.86%  '
$ws.Range("E44").Value = '  -1.60%  '
$ws.Range("E45").Value = '  +2.53%  '
$ws.Range("E46").Value = '  +4.33%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '26.62'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +3.56%  '
$ws.Range("E48").Value = '  +0.03%  '
$ws.Range("E49").Value = '  +2.44%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '35.64'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +12.29%  '
$ws.Range("E51").Value = '  +1.94%  '
